$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the preproduction environment URL/host (remove the leading "i-" prefix)
$ws.Range("A3").Value = "preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"

# Update the active selection on the sheet
$ws.Range("B4").Select()
